$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# Clear the AutoFilter criteria on the ITI sheet: this removes the
# <filterColumn> entries, unhides the rows that were hidden by the filter,
# and keeps the autoFilter range + sortState intact (same as Data > Clear).
$ws.ShowAllData()

# The edited workbook was last saved with "ITI" as the active sheet
# (previously it was "SPN") and cell E19 selected on it.
$ws.Activate()
$ws.Range("E19").Select()
